$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Mock")

# Narrow column B from 40 to 27
# (ColumnWidth uses character units that get offset by the default font's
# padding when round-tripped through the OOXML <col width> attribute, so we
# compensate to land exactly on 27.00 in the saved file)
$ws.Columns.Item(2).ColumnWidth = 26.166666666666668

# Strip the redundant group-name prefix from per-trigger option labels
# Custom: Basic Kill
$ws.Range("B40").Value = "Chance"
$ws.Range("B41").Value = "Time Scale"
$ws.Range("B42").Value = "Duration"
$ws.Range("B43").Value = "Cooldown"
$ws.Range("B44").Value = "Smoothing"
$ws.Range("B45").Value = "Third Person Distribution"

# Custom: Critical Kill
$ws.Range("B48").Value = "Chance"
$ws.Range("B49").Value = "Time Scale"
$ws.Range("B50").Value = "Duration"
$ws.Range("B51").Value = "Cooldown"
$ws.Range("B52").Value = "Smoothing"
$ws.Range("B53").Value = "Third Person Distribution"

# Custom: Dismemberment
$ws.Range("B56").Value = "Chance"
$ws.Range("B57").Value = "Time Scale"
$ws.Range("B58").Value = "Duration"
$ws.Range("B59").Value = "Cooldown"
$ws.Range("B60").Value = "Smoothing"
$ws.Range("B61").Value = "Third Person Distribution"

# Custom: Decapitation
$ws.Range("B64").Value = "Chance"
$ws.Range("B65").Value = "Time Scale"
$ws.Range("B66").Value = "Duration"
$ws.Range("B67").Value = "Cooldown"
$ws.Range("B68").Value = "Smoothing"
$ws.Range("B69").Value = "Third Person Distribution"

# Custom: Last Enemy
$ws.Range("B72").Value = "Chance"
$ws.Range("B73").Value = "Time Scale"
$ws.Range("B74").Value = "Duration"
$ws.Range("B75").Value = "Cooldown"
$ws.Range("B76").Value = "Smoothing"
$ws.Range("B77").Value = "Third Person Distribution"

# Custom: Last Stand
$ws.Range("B80").Value = "Time Scale"
$ws.Range("B81").Value = "Duration"
$ws.Range("B82").Value = "Cooldown"
$ws.Range("B83").Value = "Smoothing"

# Custom: Parry
$ws.Range("B86").Value = "Chance"
$ws.Range("B87").Value = "Time Scale"
$ws.Range("B88").Value = "Duration"
$ws.Range("B89").Value = "Cooldown"
$ws.Range("B90").Value = "Smoothing"
